$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 12
$wsExpo.Range("F3").Value = 151
$wsExpo.Range("F4").Value = 24
$wsExpo.Range("F6").Value = 4959
$wsExpo.Range("F7").Value = 4959
$wsExpo.Range("F8").Value = 50
$wsExpo.Range("F10").Value = 486
$wsExpo.Range("F13").Value = 678
$wsExpo.Range("F14").Value = 4710
$wsExpo.Range("F16").Value = 191
$wsExpo.Range("F17").Value = 199
$wsExpo.Range("F19").Value = 234
$wsExpo.Range("F20").Value = 3672
$wsExpo.Range("F23").Value = 31
$wsExpo.Range("F24").Value = 3439
$wsExpo.Range("F25").Value = 159
$wsExpo.Range("F26").Value = 147
$wsExpo.Range("F28").Value = 171
$wsExpo.Range("F29").Value = 222
$wsExpo.Range("F31").Value = 98
$wsExpo.Range("F36").Value = 6056
$wsExpo.Range("F37").Value = 953
$wsExpo.Range("F38").Value = 457
$wsExpo.Range("F42").Value = 1250
$wsExpo.Range("F43").Value = 134
$wsExpo.Range("F44").Value = 584
$wsExpo.Range("F46").Value = 2116

# --- Sheet "演出" (Performance) ---
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("G7").Value = "不可售"
$wsShow.Range("F15").Value = 136
$wsShow.Range("F22").Value = 2
$wsShow.Range("F23").Value = 784

# --- Sheet "全部类型" (All types) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value = 151
$wsAll.Range("F6").Value = 24
$wsAll.Range("F8").Value = 4959
$wsAll.Range("F9").Value = 4959
$wsAll.Range("F10").Value = 50
$wsAll.Range("F14").Value = 486
$wsAll.Range("F16").Value = 678
$wsAll.Range("F17").Value = 4710
$wsAll.Range("F18").Value = 191
$wsAll.Range("F19").Value = 199
$wsAll.Range("F21").Value = 234
$wsAll.Range("F22").Value = 3672
$wsAll.Range("F23").Value = 3439
$wsAll.Range("F24").Value = 159
$wsAll.Range("F25").Value = 147
$wsAll.Range("F26").Value = 171
$wsAll.Range("F27").Value = 222
$wsAll.Range("F29").Value = 99
$wsAll.Range("F32").Value = 136
$wsAll.Range("F34").Value = 6056
$wsAll.Range("F35").Value = 953
$wsAll.Range("F36").Value = 457
$wsAll.Range("F42").Value = 1250
$wsAll.Range("F43").Value = 134
$wsAll.Range("F44").Value = 584
$wsAll.Range("F45").Value = 2116
$wsAll.Range("F48").Value = 739

